# Generate Report for Handback
# Updates localization-status workbook to reflect a failed handback
# transform for the "23d371d4-..." item.

$wb = $excel.ActiveWorkbook

# --- Status text for the 23d371d4 item, everywhere it is shown ---------
# (Overview!E3/F3 and the per-locale sheets' Status column, C3, all
# share the same underlying string - update them together so the text
# is changed in place rather than leaving an orphaned string behind.)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Handback transform failed"
$wsOverview.Range("F3").Value = "Handback transform failed"

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = "Handback transform failed"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = "Handback transform failed"

# --- zh-cn sheet: Error Detail for row 3 --------------------------------
# (ColumnWidth 39.17 round-trips to the stored OOXML width of 40, the
# same way column A's stored width of 40 reads back as 39.17.)
$wsZh.Range("P3").Value = "Handback file name: ux5xeejh.pft is different with handoff file name: 23d371d4-33e5-42db-801e-c4de3942be76.535103aad5c0d876e5e495739b38d5ebaace7286.zh-cn."
$wsZh.Columns.Item(16).ColumnWidth = 39.17

# --- de-de sheet: Error Detail for row 3 --------------------------------
$wsDe.Range("P3").Value = "Handback file name: ux5xeejh.pft is different with handoff file name: 23d371d4-33e5-42db-801e-c4de3942be76.535103aad5c0d876e5e495739b38d5ebaace7286.de-de."
$wsDe.Columns.Item(16).ColumnWidth = 39.17
